$wb = $excel.ActiveWorkbook

# Data for the two sheets (NBR -> sheet "NBR", BAR -> sheet "BAR"), rows 3..17
# Column A: 1..15, Column B: 6..20, Column C: per-sheet reaction numbers
$nbrC = @(764, 760, 758, 745, 744, 738, 775, 776, 784, 782, 784, 774, 770, 763, 0)
$barC = @(1126, 1123, 1122, 1119, 1118, 1104, 1050, 1045, 1028, 1026, 1036, 1031, 1030, 1029, 0)

$sheets = @(
    @{ Name = "NBR"; CValues = $nbrC },
    @{ Name = "BAR"; CValues = $barC }
)

foreach ($entry in $sheets) {
    $ws = $wb.Worksheets.Item($entry.Name)

    # Extend the formatted/bordered style of column A (row 2) down through row 17
    # so the new cells pick up the same cell style (s="1") already used by A2.
    $ws.Range("A2").Copy($ws.Range("A3:A17"))

    $cValues = $entry.CValues
    for ($i = 0; $i -lt 15; $i++) {
        $r = 3 + $i
        $ws.Cells.Item($r, 1).Value = $i + 1
        $ws.Cells.Item($r, 2).Value = 6 + $i
        $ws.Cells.Item($r, 3).Value = $cValues[$i]
    }
}
